$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.174.29'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '2.061.45'
$ws.Range("E3").Value = '  -1.54%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''250.32'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").Value = '''0.678'
$ws.Range("E6").Value = '  +3.02%  '
$ws.Range("D7").Value = '''59.48'
$ws.Range("E7").Value = '  +17.13%  '
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '''60.87'
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("D10").Value = '''0.381'
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("D11").Value = '''0.0801'
$ws.Range("E11").Value = '  +7.07%  '
$ws.Range("E12").Value = '  +1.39%  '
$ws.Range("D13").Value = '''15.28'
$ws.Range("E13").Value = '  +0.54%  '
$ws.Range("D14").Value = '2.363.58'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").Value = '''0.820'
$ws.Range("E15").Value = '  -1.96%  '
$ws.Range("D16").Value = '''5.35'
$ws.Range("E16").Value = '  +4.25%  '
$ws.Range("D17").Value = '2.061.95'
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("D18").Value = '37.172.62'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = '''75.16'
$ws.Range("E19").Value = '  +3.49%  '
$ws.Range("D20").Value = '0.0₃0923'
$ws.Range("E20").Value = '  +11.86%  '
$ws.Range("E21").Value = '  +9.13%  '
$ws.Range("E22").Value = '  +2.48%  '
$ws.Range("D23").Value = '''239.39'
$ws.Range("E23").Value = '  -0.58%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").Value = '''2.44'
$ws.Range("E25").Value = '  -1.78%  '
$ws.Range("D26").Value = '''171.85'
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("D28").Value = '''20.36'
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("D29").Value = '''2.01'
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("E30").Value = '  +2.75%  '
$ws.Range("E31").Value = '  +2.70%  '
$ws.Range("E32").Value = '  -2.31%  '
$ws.Range("D33").Value = '''0.0635'
$ws.Range("E33").Value = '  +4.38%  '
$ws.Range("D34").Value = '''4.41'
$ws.Range("E34").Value = '  +7.46%  '
$ws.Range("D35").Value = '''0.0882'
$ws.Range("E35").Value = '  -5.48%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").Value = '''2.28'
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("E38").Value = '  -3.03%  '
$ws.Range("E39").Value = '  +26.47%  '
$ws.Range("E40").Value = '  +1.90%  '
$ws.Range("D41").Value = '''18.50'
$ws.Range("E41").Value = '  +4.47%  '
$ws.Range("D42").Value = '''0.0226'
$ws.Range("E42").Value = '  +1.11%  '
$ws.Range("E43").Value = '  -0.20%  '
$ws.Range("D44").Value = '''4.42'
$ws.Range("E44").Value = '  +32.36%  '
$ws.Range("D45").Value = '''97.43'
$ws.Range("E45").Value = '  -0.65%  '
$ws.Range("D46").Value = '''2.79'
$ws.Range("E46").Value = '  -0.17%  '
$ws.Range("D47").Value = '''4.48'
$ws.Range("E47").Value = '  +13.33%  '
$ws.Range("D48").Value = '''2.52'
$ws.Range("E48").Value = '  +10.92%  '
$ws.Range("D49").Value = '1.305.46'
$ws.Range("E49").Value = '  -1.00%  '
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("D51").Value = '''6.93'
$ws.Range("E51").Value = '  +0.05%  '
